$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new header cells for Wins / Losses / Ties, matching the formatting of
# the existing header row (bold, centered, bordered - same format as AC1).
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

# Fill in the team record for every data row (2-38).
for ($r = 2; $r -le 38; $r++) {
    $ws.Cells.Item($r, 30).Value = 91  # AD = col 30
    $ws.Cells.Item($r, 31).Value = 71  # AE = col 31
    $ws.Cells.Item($r, 32).Value = 0   # AF = col 32
}
